$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Some new prices are plain numeric-looking strings (e.g. '1.00', '0.540').
# Typed into a General-formatted cell, Excel would normalise these into
# numbers (dropping the formatting the source data relies on), so for those
# cells we briefly mark the cell as Text, enter the value, then restore the
# cell's style to Normal so no extra formatting is left behind.
$ws.Range("D2").Value = "51.089.19"
$ws.Range("D3").Value = "2.901.33"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "370.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "3.358.72"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "2.892.53"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.940"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "51.026.71"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "0.0₃0944"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "35.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0419"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "117.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "2.037.44"
$ws.Range("D50").Value = "3.193.40"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.240"
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  +3.21%  "
$ws.Range("E6").Value = "  -7.05%  "
$ws.Range("E7").Value = "  -6.00%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -7.61%  "
$ws.Range("E10").Value = "  -6.71%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  -5.54%  "
$ws.Range("E13").Value = "  -6.40%  "
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("E15").Value = "  -6.34%  "
$ws.Range("E16").Value = "  -3.41%  "
$ws.Range("E17").Value = "  -5.99%  "
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("E19").Value = "  -6.63%  "
$ws.Range("E20").Value = "  -5.94%  "
$ws.Range("E21").Value = "  -7.72%  "
$ws.Range("E22").Value = "  -4.81%  "
$ws.Range("E23").Value = "  -3.52%  "
$ws.Range("E24").Value = "  -5.01%  "
$ws.Range("E25").Value = "  -4.73%  "
$ws.Range("E26").Value = "  +3.77%  "
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  -6.45%  "
$ws.Range("E30").Value = "  -7.74%  "
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("E32").Value = "  -6.29%  "
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("E35").Value = "  -8.42%  "
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("E39").Value = "  -6.09%  "
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("E41").Value = "  -7.59%  "
$ws.Range("E42").Value = "  -10.16%  "
$ws.Range("E43").Value = "  -6.43%  "
$ws.Range("E44").Value = "  -7.96%  "
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("E47").Value = "  -5.36%  "
$ws.Range("E49").Value = "  -9.19%  "
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("E51").Value = "  -3.17%  "

